# Daily attendance processing - 2025-11-30 14:24:37
# Normalize the "Recorded By" (column G) text on the Session Analysis Results
# sheet: move the trailing "System" token to the front for the common
# "<email>, System" pattern, and fix the "system, System" ordering for the
# backup@backdoor.com rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
